# Week 15 simulations update
# - Adds a new player row (L.Smith) to the "Rushing" sheet, pushing the
#   existing K.Smith stat line down into a new row 10.
# - Updates accumulated season stats on both the "Rushing" and
#   "Receiving" sheets to reflect the newly simulated week.

$wb = $excel.ActiveWorkbook
$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# ---------------------------------------------------------------
# Rushing sheet updates
# ---------------------------------------------------------------

# M.Ryan (row 2)
$rushing.Range("E2").Value = 10
$rushing.Range("F2").Value = 3

# M.Davis (row 4)
$rushing.Range("C4").Value = 76
$rushing.Range("D4").Value = 28
$rushing.Range("E4").Value = 7

# C.Patterson (row 5)
$rushing.Range("C5").Value = 83
$rushing.Range("D5").Value = 32
$rushing.Range("F5").Value = 23

# Q.Ollison (row 7)
$rushing.Range("C7").Value = 12
$rushing.Range("D7").Value = 2
$rushing.Range("F7").Value = 2

# Row 9 becomes the new "L.Smith" entry (string stays same shared cell,
# but the underlying name now refers to the newly inserted player); its
# rushing stats for the week are essentially nil besides one attempt.
$rushing.Range("B9").Value = "L.Smith"
$rushing.Range("D9").Value = 0
$rushing.Range("E9").Value = 0
$rushing.Range("F9").Value = 0

# New row 10: K.Smith's stat line, carried down from the old row 9.
$rushing.Range("A9").Copy($rushing.Range("A10"))
$rushing.Range("A10").Value = 8
$rushing.Range("B10").Value = "K.Smith"
$rushing.Range("C10").Value = 1
$rushing.Range("D10").Value = 1
$rushing.Range("E10").Value = 1
$rushing.Range("F10").Value = 1

# ---------------------------------------------------------------
# Receiving sheet updates
# ---------------------------------------------------------------

# M.Davis (row 2)
$receiving.Range("C2").Value = 49
$receiving.Range("D2").Value = 37
$receiving.Range("E2").Value = 2
$receiving.Range("F2").Value = 2
$receiving.Range("G2").Value = 5

# C.Patterson (row 3)
$receiving.Range("C3").Value = 52
$receiving.Range("D3").Value = 38
$receiving.Range("G3").Value = 11

# R.Gage (row 6)
$receiving.Range("C6").Value = 60
$receiving.Range("D6").Value = 52
$receiving.Range("E6").Value = 11
$receiving.Range("F6").Value = 7

# T.Sharpe (row 9)
$receiving.Range("C9").Value = 28
$receiving.Range("D9").Value = 24
$receiving.Range("E9").Value = 5

# K.Pitts (row 10)
$receiving.Range("C10").Value = 64
$receiving.Range("D10").Value = 41
$receiving.Range("E10").Value = 23
$receiving.Range("F10").Value = 13

# H.Hurst (row 11)
$receiving.Range("C11").Value = 25
$receiving.Range("D11").Value = 21
$receiving.Range("G11").Value = 6
$receiving.Range("H11").Value = 5

# L.Smith (row 12)
$receiving.Range("C12").Value = 9
$receiving.Range("D12").Value = 8
